# Append Prithvi Shaw's remaining innings (rows 15-27) to the existing
# "Delhi Capitals" sheet, mirroring the data already present in rows 2-14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0

$newRows = @(
    @(' Dubai (DSC)', ' October 14 2020', 'Capitals won by 13 runs', 'Delhi Capitals', 'Rajasthan Royals', "Prithvi Shaw$nbsp", '0', '1', '0', '0', '0.00'),
    @(' Abu Dhabi', ' October 11 2020', 'Mumbai won by 5 wickets (with 2 balls remaining)', 'Delhi Capitals', 'Mumbai Indians', "Prithvi Shaw$nbsp", '4', '3', '1', '0', '133.33'),
    @(' Abu Dhabi', ' November 02 2020', 'Capitals won by 6 wickets (with 6 balls remaining)', 'Delhi Capitals', 'Royal Challengers Bangalore', "Prithvi Shaw$nbsp", '9', '6', '2', '0', '150.00'),
    @(' Dubai (DSC)', ' November 05 2020', 'Mumbai won by 57 runs', 'Delhi Capitals', 'Mumbai Indians', "Prithvi Shaw$nbsp", '0', '2', '0', '0', '0.00'),
    @(' Abu Dhabi', ' September 29 2020', 'Sunrisers won by 15 runs', 'Delhi Capitals', 'Sunrisers Hyderabad', "Prithvi Shaw$nbsp", '2', '5', '0', '0', '40.00'),
    @(' Sharjah', ' October 17 2020', 'Capitals won by 5 wickets (with 1 ball remaining)', 'Delhi Capitals', 'Chennai Super Kings', "Prithvi Shaw$nbsp", '0', '2', '0', '0', '0.00'),
    @(' Dubai (DSC)', ' October 31 2020', 'Mumbai won by 9 wickets (with 34 balls remaining)', 'Delhi Capitals', 'Mumbai Indians', "Prithvi Shaw$nbsp", '10', '11', '2', '0', '90.90'),
    @(' Dubai (DSC)', ' October 20 2020', 'Kings XI won by 5 wickets (with 6 balls remaining)', 'Delhi Capitals', 'Kings XI Punjab', "Prithvi Shaw$nbsp", '7', '11', '1', '0', '63.63'),
    @(' Dubai (DSC)', ' October 05 2020', 'Capitals won by 59 runs', 'Delhi Capitals', 'Royal Challengers Bangalore', "Prithvi Shaw$nbsp", '42', '23', '5', '2', '182.60'),
    @(' Dubai (DSC)', ' September 20 2020', 'Match tied (Capitals won the one-over eliminator)', 'Delhi Capitals', 'Kings XI Punjab', "Prithvi Shaw$nbsp", '5', '9', '1', '0', '55.55'),
    @(' Sharjah', ' October 03 2020', 'Capitals won by 18 runs', 'Delhi Capitals', 'Kolkata Knight Riders', "Prithvi Shaw$nbsp", '66', '41', '4', '4', '160.97'),
    @(' Dubai (DSC)', ' September 25 2020', 'Capitals won by 44 runs', 'Delhi Capitals', 'Chennai Super Kings', "Prithvi Shaw$nbsp", '64', '43', '9', '1', '148.83'),
    @(' Sharjah', ' October 09 2020', 'Capitals won by 46 runs', 'Delhi Capitals', 'Rajasthan Royals', "Prithvi Shaw$nbsp", '19', '10', '2', '1', '190.00')
)

$startRow = 15
$cols = @('A', 'B', 'C', 'D', 'E', 'F', 'G', 'H', 'I', 'J', 'K')

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    # Force text storage (matches the existing rows' numberStoredAsText
    # cells) so numeric-looking strings like "0.00" / "133.33" stay text.
    $ws.Range("A" + $r + ":K" + $r).NumberFormat = "@"

    for ($c = 0; $c -lt $cols.Count; $c++) {
        $ws.Range($cols[$c] + $r).Value = $rowData[$c]
    }
}
